$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column E (shifts old E -> F), carrying formatting from D (style 2)
$ws.Columns("E:E").Insert() | Out-Null

# Populate the new column E (Funcion/Procedimiento) for rows 1-11
$ws.Range("E1").Value = "Funcion/Procedimiento"
$ws.Range("E2").Value = "proceso.val_num_1"
$ws.Range("E3").Value = "proceso.val_num_2"
$ws.Range("E4").Value = "proceso.val_num_3"
$ws.Range("E5").Value = "proceso.val_num_4"
$ws.Range("E6").Value = "proceso.val_codigo_tabla"
$ws.Range("E7").Value = "proceso.val_codigo_tabla"
$ws.Range("E8").Value = "proceso.val_codigo_tabla"
$ws.Range("E9").Value = "proceso.val_codigo_tabla"
$ws.Range("E10").Value = "proceso.val_num_9"
$ws.Range("E11").Value = "proceso.val_num_10"

# Rows 12-14 have no function/procedure mapped yet - clear the leftover cells
# so column E stays genuinely empty there (matching the source workbook)
$ws.Range("E12:E14").Clear() | Out-Null

# Resize columns to match the widened layout
# (ColumnWidth is quantized internally to 1/6-character steps, so these are
# the closest achievable values to the author's 76.5546875 / 9.88671875 /
# 11 / 29.44140625 stored widths)
$ws.Columns("B").ColumnWidth = 75.66666666666667
$ws.Columns("C").ColumnWidth = 9
$ws.Columns("D").ColumnWidth = 10.166666666666666
$ws.Columns("E").ColumnWidth = 28.666666666666668

# Update the selection to the cell the author left active
$ws.Activate() | Out-Null
$ws.Range("E12").Select() | Out-Null
